$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 (sheet1) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("B2").Value = "2024-07-30"
$ws1.Range("C2").Value = "巢湖·元气动漫游戏嘉年华"
$ws1.Range("D2").Value = "团结东路7号 巢湖宾馆"
$ws1.Range("E2").Value = "2024.07.30 10:00-07.30 17:00"
$ws1.Range("F2").Value = 55
$ws1.Range("G2").Value = 45
$ws1.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=88193"
$ws1.Range("I2").Value = "//i2.hdslb.com/bfs/openplatform/202406/3VBeQfqQ1719318873395.jpeg"

$ws1.Range("B3").Value = "2024-08-01"
$ws1.Range("C3").Value = "合肥·ACGN夏日游园会预热场"
$ws1.Range("D3").Value = "五里墩街道长江西路与金牛路交叉口向北300米 水善汇都市微度假"
$ws1.Range("E3").Value = "2024.08.01 09:30-08.02 18:00"
$ws1.Range("F3").Value = 79
$ws1.Range("G3").Value = 45
$ws1.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=89914"
$ws1.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202407/rfRl5Bbj1721815713827.jpeg"

$ws1.Range("B4").Value = "2024-08-03"
$ws1.Range("C4").Value = "合肥·第七届环形宇宙动漫游戏嘉年华"
$ws1.Range("D4").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws1.Range("E4").Value = "2024.08.03 09:30-08.04 17:00"
$ws1.Range("F4").Value = 7188
$ws1.Range("G4").Value = 44.1
$ws1.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=84767"
$ws1.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202404/nBGuQecO1713856894035.jpeg"

$ws1.Range("B5").Value = "2024-08-10"
$ws1.Range("C5").Value = "合肥·排球少年only之夏日招新季"
$ws1.Range("D5").Value = "广德路与长江东路交口往北200米文一时埠里文旅街区 巅峰篮球公园"
$ws1.Range("E5").Value = "2024.08.10 10:00-08.10 17:00"
$ws1.Range("F5").Value = 262
$ws1.Range("G5").Value = 70
$ws1.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=88281"
$ws1.Range("I5").Value = "//i0.hdslb.com/bfs/openplatform/202406/qjd7yzXE1719556597555.jpeg"

$ws1.Range("B6").Value = "2024-08-10"
$ws1.Range("C6").Value = "合肥·比翼连枝国乙&代号鸢only"
$ws1.Range("D6").Value = "长江东大街与东二环路交叉口向南300米东方摩域商业广场三楼 格律诗婚礼艺术中心(筑梦店)"
$ws1.Range("E6").Value = "2024.08.10 09:00-08.10 22:00"
$ws1.Range("F6").Value = 412
$ws1.Range("G6").Value = 65
$ws1.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=88421"
$ws1.Range("I6").Value = "//i0.hdslb.com/bfs/openplatform/202407/RHiXT98J1721199172046.jpeg"

$ws1.Range("B7").Value = "2024-08-17"
$ws1.Range("C7").Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus"
$ws1.Range("D7").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws1.Range("E7").Value = "2024.08.17 09:30-08.18 17:00"
$ws1.Range("F7").Value = 3659
$ws1.Range("G7").Value = 69
$ws1.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=88650"
$ws1.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202407/4I7mduRV1720071650216.jpeg"

$ws1.Range("B8").Value = "2024-08-17"
$ws1.Range("C8").Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus~水千丞周边预约票"
$ws1.Range("D8").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws1.Range("E8").Value = "2024.08.17 09:30-08.17 17:00"
$ws1.Range("F8").Value = 304
$ws1.Range("G8").Value = 0.1
$ws1.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=89420"
$ws1.Range("I8").Value = "//i0.hdslb.com/bfs/openplatform/202407/hsiXAged1721203655434.jpeg"

$ws1.Range("B9").Value = "2024-08-17"
$ws1.Range("C9").Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus~水千丞签售预约票"
$ws1.Range("D9").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws1.Range("E9").Value = "2024.08.17 09:30-08.17 17:00"
$ws1.Range("F9").Value = 527
$ws1.Range("G9").Value = 0.1
$ws1.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=89421"
$ws1.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202407/r8wJqvVO1721202573195.jpeg"

$ws1.Range("B10").Value = "2024-08-17"
$ws1.Range("C10").Value = "合肥·银魂主题派对only2.0"
$ws1.Range("D10").Value = "长江东路1137号圣大国际商贸中心2-301室 梦田音乐LiveHouse(合肥店)"
$ws1.Range("E10").Value = "2024.08.17 13:00-08.17 18:00"
$ws1.Range("F10").Value = 267
$ws1.Range("G10").Value = 128
$ws1.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=87173"
$ws1.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202406/aSc8SoTl1718078234193.png"

$ws1.Range("B11").Value = "2024-08-18"
$ws1.Range("C11").Value = "合肥·SSS第五人格only"
$ws1.Range("D11").Value = "桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间"
$ws1.Range("E11").Value = "2024.08.18 09:00-08.18 17:00"
$ws1.Range("F11").Value = 597
$ws1.Range("G11").Value = 68
$ws1.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=88430"
$ws1.Range("I11").Value = "//i0.hdslb.com/bfs/openplatform/202406/a0qh8I1h1719660853555.png"

$ws1.Range("B12").Value = "2024-09-07"
$ws1.Range("C12").Value = "合肥·国乙only宇宙心动（含夜场）"
$ws1.Range("D12").Value = "文忠路1865号 赫拉诺言艺术中心"
$ws1.Range("E12").Value = "2024.09.07 10:00-09.07 21:00"
$ws1.Range("F12").Value = 88
$ws1.Range("G12").Value = 48
$ws1.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=89803"
$ws1.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202407/w5hQDj821721564303601.jpeg"

# delete now-unused trailing rows 13:16 in sheet1 (shrinks dimension to A1:I12)
$ws1.Range("A13:I16").Clear()

# ---- Sheet: 全部类型 (sheet4) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("B2").Value = "2024-07-26"
$ws4.Range("C2").Value = "合肥·Yolo Fes永乐庆典Vol.3·少女偶像联合演出DAY1&DAY3"
$ws4.Range("D2").Value = "金寨路与天堂窄路交叉口 梵木艺术中心"
$ws4.Range("E2").Value = "2024.07.26 18:00-07.28 23:59"
$ws4.Range("F2").Value = 12
$ws4.Range("G2").Value = 178
$ws4.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=89514"
$ws4.Range("I2").Value = "//i2.hdslb.com/bfs/openplatform/202407/aMtLMGR31721289854139.jpeg"

$ws4.Range("B3").Value = "2024-07-30"
$ws4.Range("C3").Value = "巢湖·元气动漫游戏嘉年华"
$ws4.Range("D3").Value = "团结东路7号 巢湖宾馆"
$ws4.Range("E3").Value = "2024.07.30 10:00-07.30 17:00"
$ws4.Range("F3").Value = 55
$ws4.Range("G3").Value = 45
$ws4.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=88193"
$ws4.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202406/3VBeQfqQ1719318873395.jpeg"

$ws4.Range("B4").Value = "2024-08-01"
$ws4.Range("C4").Value = "合肥·ACGN夏日游园会预热场"
$ws4.Range("D4").Value = "五里墩街道长江西路与金牛路交叉口向北300米 水善汇都市微度假"
$ws4.Range("E4").Value = "2024.08.01 09:30-08.02 18:00"
$ws4.Range("F4").Value = 79
$ws4.Range("G4").Value = 45
$ws4.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=89914"
$ws4.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202407/rfRl5Bbj1721815713827.jpeg"

$ws4.Range("B5").Value = "2024-08-02"
$ws4.Range("C5").Value = "合肥·新西兰·治愈系民谣歌手Luke Thompson2024中国巡演 KEEP ROLLING ON "
$ws4.Range("D5").Value = "宁国路罍街二期15号楼安徽原创音乐基地3楼 合肥ON THE WAY LiveHouse"
$ws4.Range("E5").Value = "2024.08.02 20:00-08.02 21:30"
$ws4.Range("F5").Value = 4
$ws4.Range("G5").Value = 180
$ws4.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=88824"
$ws4.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202407/FKItJRNl1719803666645.jpeg"

$ws4.Range("B6").Value = "2024-08-03"
$ws4.Range("C6").Value = "合肥·第七届环形宇宙动漫游戏嘉年华"
$ws4.Range("D6").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws4.Range("E6").Value = "2024.08.03 09:30-08.04 17:00"
$ws4.Range("F6").Value = 7188
$ws4.Range("G6").Value = 44.1
$ws4.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=84767"
$ws4.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202404/nBGuQecO1713856894035.jpeg"

$ws4.Range("B7").Value = "2024-08-03"
$ws4.Range("C7").Value = "合肥·首届包河留声机音乐节—《菊次郎的夏天》久石让钢琴曲梦幻之旅演奏会"
$ws4.Range("D7").Value = "徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院"
$ws4.Range("E7").Value = "2024.08.03 19:30-08.03 21:00"
$ws4.Range("F7").Value = 49
$ws4.Range("G7").Value = 80
$ws4.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=83556"
$ws4.Range("I7").Value = "//i1.hdslb.com/bfs/openplatform/202403/4nwOTVDu1711695345941.jpeg"

$ws4.Range("B8").Value = "2024-08-10"
$ws4.Range("C8").Value = "合肥·排球少年only之夏日招新季"
$ws4.Range("D8").Value = "广德路与长江东路交口往北200米文一时埠里文旅街区 巅峰篮球公园"
$ws4.Range("E8").Value = "2024.08.10 10:00-08.10 17:00"
$ws4.Range("F8").Value = 262
$ws4.Range("G8").Value = 70
$ws4.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=88281"
$ws4.Range("I8").Value = "//i0.hdslb.com/bfs/openplatform/202406/qjd7yzXE1719556597555.jpeg"

$ws4.Range("B9").Value = "2024-08-10"
$ws4.Range("C9").Value = "合肥·比翼连枝国乙&代号鸢only"
$ws4.Range("D9").Value = "长江东大街与东二环路交叉口向南300米东方摩域商业广场三楼 格律诗婚礼艺术中心(筑梦店)"
$ws4.Range("E9").Value = "2024.08.10 09:00-08.10 22:00"
$ws4.Range("F9").Value = 412
$ws4.Range("G9").Value = 65
$ws4.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=88421"
$ws4.Range("I9").Value = "//i0.hdslb.com/bfs/openplatform/202407/RHiXT98J1721199172046.jpeg"

$ws4.Range("B10").Value = "2024-08-17"
$ws4.Range("C10").Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus"
$ws4.Range("D10").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws4.Range("E10").Value = "2024.08.17 09:30-08.18 17:00"
$ws4.Range("F10").Value = 3659
$ws4.Range("G10").Value = 69
$ws4.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=88650"
$ws4.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202407/4I7mduRV1720071650216.jpeg"

$ws4.Range("B11").Value = "2024-08-17"
$ws4.Range("C11").Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus~水千丞周边预约票"
$ws4.Range("D11").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws4.Range("E11").Value = "2024.08.17 09:30-08.17 17:00"
$ws4.Range("F11").Value = 304
$ws4.Range("G11").Value = 0.1
$ws4.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=89420"
$ws4.Range("I11").Value = "//i0.hdslb.com/bfs/openplatform/202407/hsiXAged1721203655434.jpeg"

$ws4.Range("B12").Value = "2024-08-17"
$ws4.Range("C12").Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus~水千丞签售预约票"
$ws4.Range("D12").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws4.Range("E12").Value = "2024.08.17 09:30-08.17 17:00"
$ws4.Range("F12").Value = 527
$ws4.Range("G12").Value = 0.1
$ws4.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=89421"
$ws4.Range("I12").Value = "//i2.hdslb.com/bfs/openplatform/202407/r8wJqvVO1721202573195.jpeg"

$ws4.Range("B13").Value = "2024-08-17"
$ws4.Range("C13").Value = "合肥·银魂主题派对only2.0"
$ws4.Range("D13").Value = "长江东路1137号圣大国际商贸中心2-301室 梦田音乐LiveHouse(合肥店)"
$ws4.Range("E13").Value = "2024.08.17 13:00-08.17 18:00"
$ws4.Range("F13").Value = 267
$ws4.Range("G13").Value = 128
$ws4.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=87173"
$ws4.Range("I13").Value = "//i2.hdslb.com/bfs/openplatform/202406/aSc8SoTl1718078234193.png"

$ws4.Range("B14").Value = "2024-08-18"
$ws4.Range("C14").Value = "合肥·SSS第五人格only"
$ws4.Range("D14").Value = "桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间"
$ws4.Range("E14").Value = "2024.08.18 09:00-08.18 17:00"
$ws4.Range("F14").Value = 597
$ws4.Range("G14").Value = 68
$ws4.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=88430"
$ws4.Range("I14").Value = "//i0.hdslb.com/bfs/openplatform/202406/a0qh8I1h1719660853555.png"

$ws4.Range("B15").Value = "2024-09-07"
$ws4.Range("C15").Value = "合肥·国乙only宇宙心动（含夜场）"
$ws4.Range("D15").Value = "文忠路1865号 赫拉诺言艺术中心"
$ws4.Range("E15").Value = "2024.09.07 10:00-09.07 21:00"
$ws4.Range("F15").Value = 88
$ws4.Range("G15").Value = 48
$ws4.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=89803"
$ws4.Range("I15").Value = "//i1.hdslb.com/bfs/openplatform/202407/w5hQDj821721564303601.jpeg"

# delete now-unused trailing rows 16:19 in sheet4 (shrinks dimension to A1:I15)
$ws4.Range("A16:I19").Clear()
